$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "No utente reale (usabilità/accessibilità/" ->
#           two runs: "No utente reale (usabilità/accessibilità" + ")"
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("No utente reale (usabilità/accessibilità/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $start1 = $r1.Start
    $r1.Text = ""
    $anchor1 = $d.Range($start1, $start1)
    $anchor1.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>No utente reale (usabilit&#224;/accessibilit&#224;</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
}

# ---------------------------------------------------------------------------
# Change 2: merge "15 – 10 : " + "Starting" + " meeting" (3 runs, with
# spell-check proofErr markers around "Starting") into a single run reading
# "15 – 10 : Starting meeting", dropping the proofErr elements.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("15 – 10 : Starting meeting", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $startPoint = $d.Range(0, 0)
    $startPoint.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="38F289B3" w14:textId="594357A7" w:rsidR="001760F3" w:rsidRPr="001760F3" w:rsidRDefault="001760F3" w:rsidP="001760F3"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="001760F3"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>15 &#8211; 10 : Starting meeting</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
    $oldPara = $d.Paragraphs(2)
    $oldPara.Range.Delete()
}
